# Update countries & provincias Spain
# Applies the 13-May-2020 22:35 data refresh to the "Pais" sheet:
#  - updates the "datos actualizados" timestamp in the title cell
#  - refreshes numeric stats for several countries
#  - re-sorts a few country pairs whose case counts changed rank
#    (Peru/Canada, Sudafrica/Filipinas, Nueva Caledonia/Belice)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Title / timestamp row -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 13 de Mayo de 2020 a las 22:35"

# --- Estados Unidos (row 4) -------------------------------------------------
$ws.Range("B4").Value = 1423641
$ws.Range("C4").Value = 15005
$ws.Range("E4").Value = 1032596
$ws.Range("F4").Value = 16349
$ws.Range("G4").Value = 1322
$ws.Range("H4").Value = 84747

# --- Brasil (row 9) ----------------------------------------------------------
$ws.Range("B9").Value = 181518
$ws.Range("C9").Value = 3916
$ws.Range("E9").Value = 96218
$ws.Range("G9").Value = 299
$ws.Range("H9").Value = 12703

# --- Peru overtakes Canada (rows 16/17) -------------------------------------
# Row 16 becomes Peru with its refreshed totals
$ws.Range("A16").Value = "Peru"
$ws.Range("B16").Value = 76306
$ws.Range("C16").Value = 4247
$ws.Range("D16").Value = 24324
$ws.Range("E16").Value = 49813
$ws.Range("F16").Value = 806
$ws.Range("G16").Value = 112
$ws.Range("H16").Value = 2169

# Row 17 becomes Canada, keeping its previous (unchanged) totals
$ws.Range("A17").Value = "Canada"
$ws.Range("B17").Value = 72196
$ws.Range("C17").Value = 1039
$ws.Range("D17").Value = 34916
$ws.Range("E17").Value = 31979
$ws.Range("F17").Value = 502
$ws.Range("G17").Value = 132
$ws.Range("H17").Value = 5301

# --- Irlanda (row 31) --------------------------------------------------------
$ws.Range("F31").Value = 69

# --- Sudafrica overtakes Filipinas (rows 42/43) -----------------------------
# Row 42 becomes Sudafrica with its refreshed totals
$ws.Range("A42").Value = "Sudafrica"
$ws.Range("B42").Value = 12047
$ws.Range("C42").Value = 697
$ws.Range("D42").Value = 4745
$ws.Range("E42").Value = 7083
$ws.Range("F42").Value = 119
$ws.Range("G42").Value = 13
$ws.Range("H42").Value = 219

# Row 43 becomes Filipinas, keeping its previous (unchanged) totals
$ws.Range("A43").Value = "Filipinas"
$ws.Range("B43").Value = 11618
$ws.Range("C43").Value = 268
$ws.Range("D43").Value = 2251
$ws.Range("E43").Value = 8595
$ws.Range("F43").Value = 77
$ws.Range("G43").Value = 21
$ws.Range("H43").Value = 772

# --- Islas Caimanes (row 167) ------------------------------------------------
$ws.Range("B167").Value = 86
$ws.Range("C167").Value = 1
$ws.Range("E167").Value = 35

# --- Nueva Caledonia overtakes Belice (rows 193/194) ------------------------
# Row 193 becomes Nueva Caledonia, keeping its previous (unchanged) totals
$ws.Range("A193").Value = "Nueva Caledonia"
$ws.Range("B193").Value = 18
$ws.Range("C193").Value = 0
$ws.Range("D193").Value = 18
$ws.Range("E193").Value = 0
$ws.Range("F193").Value = 0
$ws.Range("G193").Value = 0
$ws.Range("H193").Value = 0

# Row 194 becomes Belice, keeping its previous (unchanged) totals
$ws.Range("A194").Value = "Belice"
$ws.Range("B194").Value = 18
$ws.Range("C194").Value = 0
$ws.Range("D194").Value = 16
$ws.Range("E194").Value = 0
$ws.Range("F194").Value = 0
$ws.Range("G194").Value = 0
$ws.Range("H194").Value = 2
